$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Portfolio Name" header in column L (shared string #11)
$ws.Range("L1").Value = "Portfolio Name"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").Font.Underline = $true
$ws.Range("L1").HorizontalAlignment = -4131
$ws.Range("L1").VerticalAlignment = -4108
$ws.Range("L1").WrapText = $true

# Move the active selection to M2, matching the saved view state
$ws.Range("M2").Select()
